# Add a "Brand" column to the Partner invoice Annexure template.
# A new column is inserted at column I (shifting Category/Size/City/... one
# column to the right), with a "Brand" header (row 16) and the
# "{booking:appliance_brand}" merge placeholder (row 17) underneath it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I; this shifts existing columns I:S to J:T and
# carries along column widths / cell styles the same way Excel's own
# "Insert Column" command does.
$ws.Columns("I:I").Insert()

# Populate the new column's header + templated placeholder row.
$ws.Range("I16").Value = "Brand"
$ws.Range("I17").Value = "{booking:appliance_brand}"
